# Insert a new record row at row 254 (shifting existing rows 254:311 down to 255:312)
# and populate it with the new Ciboulette price-record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(254).Insert()

$ws.Cells.Item(254, 1).Value = 4
$ws.Cells.Item(254, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(254, 3).Value = "Los Lagos"
$ws.Cells.Item(254, 4).Value = 44964
$ws.Cells.Item(254, 5).Value = 10
$ws.Cells.Item(254, 6).Value = 100112039
$ws.Cells.Item(254, 7).Value = "Ciboulette"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 240
$ws.Cells.Item(254, 11).Value = 3500
$ws.Cells.Item(254, 12).Value = 3500
$ws.Cells.Item(254, 13).Value = 3500
$ws.Cells.Item(254, 14).Value = "$/docena de atados"
$ws.Cells.Item(254, 15).Value = "Región Metropolitana"
$ws.Cells.Item(254, 16).Value = 1167
$ws.Cells.Item(254, 17).Value = 3
$ws.Cells.Item(254, 18).Value = "Hortaliza"

# Match the date-style formatting used by the rest of column D
$ws.Cells.Item(254, 4).NumberFormat = $ws.Cells.Item(255, 4).NumberFormat
